$d = $word.ActiveDocument

# Bump version string "6.26" -> "6.27" wherever it appears as a contiguous
# run of text fragments: "LogZilla_SyslogAgent_" + "6" + ".2" + "6" + ...
# and "LogZilla Windows Syslog Agent V" + "6" + ".2" + "6" + ... (footers)
#
# Use Find/Replace across the whole document (body + headers/footers) so it
# works regardless of how the text is split across runs.

$d.Content.Find.Execute("LogZilla_SyslogAgent_6.26", $true, $false, $false, $false, $false, $true, 1, $false, "LogZilla_SyslogAgent_6.27", 2)

$d.Content.Find.Execute("LogZilla Windows Syslog Agent V6.26", $true, $false, $false, $false, $false, $true, 1, $false, "LogZilla Windows Syslog Agent V6.27", 2)

foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("LogZilla Windows Syslog Agent V6.26", $true, $false, $false, $false, $false, $true, 1, $false, "LogZilla Windows Syslog Agent V6.27", 2)
        }
    }
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("LogZilla Windows Syslog Agent V6.26", $true, $false, $false, $false, $false, $true, 1, $false, "LogZilla Windows Syslog Agent V6.27", 2)
        }
    }
}
